# Insert a new data row at row 16 (pushing the existing rows 16-40 down to 17-41)
# and populate it with the new Chirimoya price-record for the week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 16..40 down by inserting a new blank row at position 16.
$ws.Rows.Item(16).Insert()

# New row 16 values (mirrors the layout/template of the surrounding rows).
$ws.Cells.Item(16, 1).Value = 11
$ws.Cells.Item(16, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(16, 3).Value = "Bíobío"
$ws.Cells.Item(16, 4).Value = 44868
$ws.Cells.Item(16, 4).NumberFormat = $ws.Cells.Item(17, 4).NumberFormat
$ws.Cells.Item(16, 5).Value = 8
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100107
$ws.Cells.Item(16, 8).Value = "Otros"
$ws.Cells.Item(16, 9).Value = 100107002
$ws.Cells.Item(16, 10).Value = "Chirimoya"
$ws.Cells.Item(16, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(16, 12).Value = "Especial"
$ws.Cells.Item(16, 13).Value = 180
$ws.Cells.Item(16, 14).Value = 24000
$ws.Cells.Item(16, 15).Value = 25000
$ws.Cells.Item(16, 16).Value = 24444
$ws.Cells.Item(16, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(16, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(16, 19).Value = 2444
$ws.Cells.Item(16, 20).Value = 10
